$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# Update header timestamps (Ultima actualizacion) on all sheets
$ws1.Cells.Item(2,1).Value = "Última actualización: 19:54:49"
$ws2.Cells.Item(2,1).Value = "Última actualización: 19:54:49"
$ws3.Cells.Item(2,1).Value = "Última actualización: 19:54:49"

# Update total row count on sheet1
$ws1.Cells.Item(3,1).Value = "Total filas: 125"

# Rewrite affected data rows on sheet1 (rows 19-20 swap, and rows 102-130 shift/insert)
$ws1.Cells.Item(19,1).Value = "16:50:41"
$ws1.Cells.Item(19,2).Value = "17:17"
$ws1.Cells.Item(19,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(19,4).Value = 27
$ws1.Cells.Item(19,5).Value = "LP1912"
$ws1.Cells.Item(20,1).Value = "16:46:42"
$ws1.Cells.Item(20,2).Value = "17:17"
$ws1.Cells.Item(20,3).Value = "17_ROMERO"
$ws1.Cells.Item(20,4).Value = 31
$ws1.Cells.Item(20,5).Value = "LP1912"
$ws1.Cells.Item(102,1).Value = "19:54:49"
$ws1.Cells.Item(102,2).Value = "19:54"
$ws1.Cells.Item(102,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(102,4).Value = 0
$ws1.Cells.Item(102,5).Value = "LP1912"
$ws1.Cells.Item(103,1).Value = "18:10:41"
$ws1.Cells.Item(103,2).Value = "19:58"
$ws1.Cells.Item(103,3).Value = "14X44_ABASTO"
$ws1.Cells.Item(103,4).Value = 108
$ws1.Cells.Item(103,5).Value = "LP1912"
$ws1.Cells.Item(104,1).Value = "18:31:18"
$ws1.Cells.Item(104,2).Value = "19:59"
$ws1.Cells.Item(104,3).Value = "14X44_ABASTO"
$ws1.Cells.Item(104,4).Value = 88
$ws1.Cells.Item(104,5).Value = "LP1912"
$ws1.Cells.Item(105,1).Value = "18:10:41"
$ws1.Cells.Item(105,2).Value = "20:00"
$ws1.Cells.Item(105,3).Value = "215C_EL PATO"
$ws1.Cells.Item(105,4).Value = 110
$ws1.Cells.Item(105,5).Value = "LP1912"
$ws1.Cells.Item(106,1).Value = "19:47:58"
$ws1.Cells.Item(106,2).Value = "20:00"
$ws1.Cells.Item(106,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(106,4).Value = 13
$ws1.Cells.Item(106,5).Value = "LP1912"
$ws1.Cells.Item(107,1).Value = "18:31:18"
$ws1.Cells.Item(107,2).Value = "20:01"
$ws1.Cells.Item(107,3).Value = "215C_EL PATO"
$ws1.Cells.Item(107,4).Value = 90
$ws1.Cells.Item(107,5).Value = "LP1912"
$ws1.Cells.Item(108,1).Value = "19:47:58"
$ws1.Cells.Item(108,2).Value = "20:02"
$ws1.Cells.Item(108,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(108,4).Value = 15
$ws1.Cells.Item(108,5).Value = "LP1912"
$ws1.Cells.Item(109,1).Value = "19:11:59"
$ws1.Cells.Item(109,2).Value = "20:04"
$ws1.Cells.Item(109,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(109,4).Value = 53
$ws1.Cells.Item(109,5).Value = "LP1912"
$ws1.Cells.Item(110,1).Value = "19:47:58"
$ws1.Cells.Item(110,2).Value = "20:09"
$ws1.Cells.Item(110,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(110,4).Value = 22
$ws1.Cells.Item(110,5).Value = "LP1912"
$ws1.Cells.Item(111,1).Value = "19:35:31"
$ws1.Cells.Item(111,2).Value = "20:10"
$ws1.Cells.Item(111,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(111,4).Value = 35
$ws1.Cells.Item(111,5).Value = "LP1912"
$ws1.Cells.Item(112,1).Value = "19:11:59"
$ws1.Cells.Item(112,2).Value = "20:13"
$ws1.Cells.Item(112,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(112,4).Value = 62
$ws1.Cells.Item(112,5).Value = "LP1912"
$ws1.Cells.Item(113,1).Value = "18:31:18"
$ws1.Cells.Item(113,2).Value = "20:14"
$ws1.Cells.Item(113,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(113,4).Value = 103
$ws1.Cells.Item(113,5).Value = "LP1912"
$ws1.Cells.Item(114,1).Value = "19:11:59"
$ws1.Cells.Item(114,2).Value = "20:25"
$ws1.Cells.Item(114,3).Value = "15_ABASTO"
$ws1.Cells.Item(114,4).Value = 74
$ws1.Cells.Item(114,5).Value = "LP1912"
$ws1.Cells.Item(115,1).Value = "18:31:18"
$ws1.Cells.Item(115,2).Value = "20:26"
$ws1.Cells.Item(115,3).Value = "15_ABASTO"
$ws1.Cells.Item(115,4).Value = 115
$ws1.Cells.Item(115,5).Value = "LP1912"
$ws1.Cells.Item(116,1).Value = "18:44:34"
$ws1.Cells.Item(116,2).Value = "20:28"
$ws1.Cells.Item(116,3).Value = "10_OLMOS"
$ws1.Cells.Item(116,4).Value = 104
$ws1.Cells.Item(116,5).Value = "LP1912"
$ws1.Cells.Item(117,1).Value = "18:31:18"
$ws1.Cells.Item(117,2).Value = "20:29"
$ws1.Cells.Item(117,3).Value = "10_OLMOS"
$ws1.Cells.Item(117,4).Value = 118
$ws1.Cells.Item(117,5).Value = "LP1912"
$ws1.Cells.Item(118,1).Value = "19:11:59"
$ws1.Cells.Item(118,2).Value = "20:43"
$ws1.Cells.Item(118,3).Value = "215B_EL PATO"
$ws1.Cells.Item(118,4).Value = 92
$ws1.Cells.Item(118,5).Value = "LP1912"
$ws1.Cells.Item(119,1).Value = "19:11:59"
$ws1.Cells.Item(119,2).Value = "20:44"
$ws1.Cells.Item(119,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(119,4).Value = 93
$ws1.Cells.Item(119,5).Value = "LP1912"
$ws1.Cells.Item(120,1).Value = "18:52:04"
$ws1.Cells.Item(120,2).Value = "20:44"
$ws1.Cells.Item(120,3).Value = "215B_EL PATO"
$ws1.Cells.Item(120,4).Value = 112
$ws1.Cells.Item(120,5).Value = "LP1912"
$ws1.Cells.Item(121,1).Value = "18:52:04"
$ws1.Cells.Item(121,2).Value = "20:45"
$ws1.Cells.Item(121,3).Value = "17X38_ROMERO"
$ws1.Cells.Item(121,4).Value = 113
$ws1.Cells.Item(121,5).Value = "LP1912"
$ws1.Cells.Item(122,1).Value = "19:54:49"
$ws1.Cells.Item(122,2).Value = "20:50"
$ws1.Cells.Item(122,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(122,4).Value = 56
$ws1.Cells.Item(122,5).Value = "LP1912"
$ws1.Cells.Item(123,1).Value = "19:35:31"
$ws1.Cells.Item(123,2).Value = "20:52"
$ws1.Cells.Item(123,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(123,4).Value = 77
$ws1.Cells.Item(123,5).Value = "LP1912"
$ws1.Cells.Item(124,1).Value = "19:54:49"
$ws1.Cells.Item(124,2).Value = "20:56"
$ws1.Cells.Item(124,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(124,4).Value = 62
$ws1.Cells.Item(124,5).Value = "LP1912"
$ws1.Cells.Item(125,1).Value = "19:11:59"
$ws1.Cells.Item(125,2).Value = "21:01"
$ws1.Cells.Item(125,3).Value = "215A_EL PATO"
$ws1.Cells.Item(125,4).Value = 110
$ws1.Cells.Item(125,5).Value = "LP1912"
$ws1.Cells.Item(126,1).Value = "19:11:59"
$ws1.Cells.Item(126,2).Value = "21:02"
$ws1.Cells.Item(126,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(126,4).Value = 111
$ws1.Cells.Item(126,5).Value = "LP1912"
$ws1.Cells.Item(127,1).Value = "19:47:58"
$ws1.Cells.Item(127,2).Value = "21:06"
$ws1.Cells.Item(127,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(127,4).Value = 79
$ws1.Cells.Item(127,5).Value = "LP1912"
$ws1.Cells.Item(128,1).Value = "19:35:31"
$ws1.Cells.Item(128,2).Value = "21:10"
$ws1.Cells.Item(128,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(128,4).Value = 95
$ws1.Cells.Item(128,5).Value = "LP1912"
$ws1.Cells.Item(129,1).Value = "19:35:31"
$ws1.Cells.Item(129,2).Value = "21:23"
$ws1.Cells.Item(129,3).Value = "10_OLMOS"
$ws1.Cells.Item(129,4).Value = 108
$ws1.Cells.Item(129,5).Value = "LP1912"
$ws1.Cells.Item(130,1).Value = "19:54:49"
$ws1.Cells.Item(130,2).Value = "21:49"
$ws1.Cells.Item(130,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(130,4).Value = 115
$ws1.Cells.Item(130,5).Value = "LP1912"
